$d = $word.ActiveDocument

$replacements = @(
    @("2025-11-14 Friday", "2025-11-15 Saturday"),
    @("753÷4=188, 1", "786÷5=157, 1"),
    @("804÷5=160, 4", "920÷7=131, 3"),
    @("541÷5=108, 1", "468÷3=156, 0"),
    @("202÷4=50, 2", "596÷4=149, 0"),
    @("811÷6=135, 1", "519÷8=64, 7"),
    @("870÷4=217, 2", "647÷6=107, 5"),
    @("802÷4=200, 2", "824÷6=137, 2"),
    @("434÷2=217, 0", "449÷9=49, 8"),
    @("967÷9=107, 4", "271÷8=33, 7"),
    @("197÷6=32, 5", "170÷9=18, 8"),
    @("575÷8=71, 7", "580÷3=193, 1"),
    @("469÷4=117, 1", "898÷7=128, 2"),
    @("564÷2=282, 0", "172÷8=21, 4"),
    @("833÷9=92, 5", "276÷9=30, 6"),
    @("848÷6=141, 2", "856÷4=214, 0"),
    @("216÷3=72, 0", "357÷6=59, 3"),
    @("234÷8=29, 2", "134÷6=22, 2"),
    @("907÷7=129, 4", "808÷2=404, 0"),
    @("254÷2=127, 0", "474÷6=79, 0"),
    @("438÷4=109, 2", "858÷7=122, 4"),
    @("749÷3=249, 2", "400÷9=44, 4"),
    @("668÷4=167, 0", "765÷8=95, 5"),
    @("586÷5=117, 1", "527÷4=131, 3"),
    @("235÷5=47, 0", "905÷4=226, 1"),
    @("629÷9=69, 8", "488÷7=69, 5")
)

foreach ($pair in $replacements) {
    $old = $pair[0]
    $new = $pair[1]
    $range = $d.Content
    $range.Find.Execute($old, $true, $false, $false, $false, $false, $true, 1, $false, $new, 2)
}
